$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "23.252.18"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  +1.00%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.604.70"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  +0.25%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9999"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  -0.11%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "1.000"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  -0.02%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "303.79"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  +0.77%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.3771"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  -0.31%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "51.98"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  +4.01%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.3634"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  -0.16%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.278"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  +1.29%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.08146"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  +0.12%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.000"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  -0.08%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "22.85"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  +1.14%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.593"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  -0.05%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.418"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  +0.86%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.00001252"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  +0.61%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "1.602.34"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  -0.23%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "94.06"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  +2.30%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06919"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  +1.43%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "18.16"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  -0.44%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.537"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  -0.29%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "1.003"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  +0.20%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "12.92"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  -1.31%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "23.245.81"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  +0.96%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.450"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  +4.04%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "3.046"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  +8.25%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "21.23"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  +0.76%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "149.81"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  -0.40%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "5.282"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  +0.65%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "135.60"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  +0.89%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "2.380"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  +3.28%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "6.762"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  -1.26%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.778.99"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  -0.34%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.9652"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  -0.35%  "

$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  -1.43%  "

$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  +0.15%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.02753"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  +1.53%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.2530"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  -0.27%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "6.147"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  -2.08%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.08804"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  -1.07%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.385"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  +1.03%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.7114"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  +1.02%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "12.49"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  +0.24%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "15.55"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  +1.79%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.6549"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  -1.54%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.321"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  +0.38%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.9993"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  -0.02%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "4.010"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  +0.49%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "132.95"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  +0.42%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.07940"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  +0.30%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.209"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  -2.15%  "

